$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3041.1724  # H64: 3065.276 -> 3041.1724
$ws.Cells.Item(64, 9).Value = 2844.111  # I64: 2866.3333 -> 2844.111
$ws.Cells.Item(64, 10).Value = 3363.6365  # J64: 3390.818 -> 3363.6365
$ws.Cells.Item(64, 11).Value = 2844.111  # K64: 2866.3333 -> 2844.111
$ws.Cells.Item(64, 12).Value = 3363.6365  # L64: 3390.818 -> 3363.6365
$ws.Cells.Item(64, 13).Value = -2596.111  # M64: -2618.3333 -> -2596.111
$ws.Cells.Item(64, 14).Value = -3859.6365  # N64: -3886.818 -> -3859.6365
$ws.Cells.Item(67, 8).Value = 3041.1724  # H67: 3065.276 -> 3041.1724
$ws.Cells.Item(67, 9).Value = 2844.111  # I67: 2866.3333 -> 2844.111
$ws.Cells.Item(67, 10).Value = 3363.6365  # J67: 3390.818 -> 3363.6365
$ws.Cells.Item(67, 11).Value = 2844.111  # K67: 2866.3333 -> 2844.111
$ws.Cells.Item(67, 12).Value = 3363.6365  # L67: 3390.818 -> 3363.6365
$ws.Cells.Item(67, 13).Value = -1986.111  # M67: -2008.3333 -> -1986.111
$ws.Cells.Item(67, 14).Value = -5079.636500000001  # N67: -5106.818 -> -5079.636500000001
$ws.Cells.Item(76, 8).Value = 25643948  # H76: 27788786 -> 25643948
$ws.Cells.Item(76, 9).Value = 3157.9  # I76: 13966.556 -> 3157.9
$ws.Cells.Item(76, 11).Value = 3157.9  # K76: 13966.556 -> 3157.9
$ws.Cells.Item(76, 13).Value = -2842.9  # M76: -13651.556 -> -2842.9
$ws.Cells.Item(79, 8).Value = 25643948  # H79: 27788786 -> 25643948
$ws.Cells.Item(79, 9).Value = 3157.9  # I79: 13966.556 -> 3157.9
$ws.Cells.Item(79, 11).Value = 3157.9  # K79: 13966.556 -> 3157.9
$ws.Cells.Item(79, 13).Value = -2065.9  # M79: -12874.556 -> -2065.9
$ws.Cells.Item(123, 8).Value = 76127.21000000001  # H123: 17639.5 -> 76127.21000000001
$ws.Cells.Item(123, 10).Value = 76127.21000000001  # J123: 17639.5 -> 76127.21000000001
$ws.Cells.Item(123, 12).Value = 76127.21000000001  # L123: 17639.5 -> 76127.21000000001
$ws.Cells.Item(123, 14).Value = -85927.21000000001  # N123: -27439.5 -> -85927.21000000001
$ws.Cells.Item(124, 8).Value = 49781  # H124: 0 -> 49781
$ws.Cells.Item(124, 10).Value = 49781  # J124: 0 -> 49781
$ws.Cells.Item(124, 12).Value = 49781  # L124: 0 -> 49781
$ws.Cells.Item(124, 14).Value = -59601  # N124: None -> -59601
$ws.Cells.Item(138, 8).Value = 3005.4353  # H138: 3129.85 -> 3005.4353
$ws.Cells.Item(138, 9).Value = 1527.7576  # I138: 1642.8966 -> 1527.7576
$ws.Cells.Item(138, 10).Value = 3943.1924  # J138: 3975.3726 -> 3943.1924
$ws.Cells.Item(138, 11).Value = 4583.2728  # K138: 4928.6898 -> 4583.2728
$ws.Cells.Item(138, 12).Value = 11829.5772  # L138: 11926.1178 -> 11829.5772
$ws.Cells.Item(138, 13).Value = 556.7272000000003  # M138: 211.3101999999999 -> 556.7272000000003
$ws.Cells.Item(138, 14).Value = -22109.5772  # N138: -22206.1178 -> -22109.5772
$ws.Cells.Item(140, 8).Value = 53735.293  # H140: 59923.53 -> 53735.293
$ws.Cells.Item(140, 10).Value = 53735.293  # J140: 59923.53 -> 53735.293
$ws.Cells.Item(140, 12).Value = 53735.293  # L140: 59923.53 -> 53735.293
$ws.Cells.Item(140, 14).Value = -64095.293  # N140: -70283.53 -> -64095.293

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 255.66667  # H4: 450 -> 255.66667
$ws.Cells.Item(4, 9).Value = 220.14285  # I4: 0 -> 220.14285
$ws.Cells.Item(4, 10).Value = 380  # J4: 450 -> 380
$ws.Cells.Item(4, 11).Value = 220.14285  # K4: 0 -> 220.14285
$ws.Cells.Item(4, 12).Value = 380  # L4: 450 -> 380
$ws.Cells.Item(4, 13).Value = -104.14285  # M4: None -> -104.14285
$ws.Cells.Item(4, 14).Value = -612  # N4: -682 -> -612
$ws.Cells.Item(74, 8).Value = 7466419.5  # H74: 11365139 -> 7466419.5
$ws.Cells.Item(74, 9).Value = 9806548  # I74: 14707484 -> 9806548
$ws.Cells.Item(74, 10).Value = 7259.75  # J74: 1164.2 -> 7259.75
$ws.Cells.Item(74, 11).Value = 9806548  # K74: 14707484 -> 9806548
$ws.Cells.Item(74, 12).Value = 7259.75  # L74: 1164.2 -> 7259.75
$ws.Cells.Item(74, 13).Value = -9805674  # M74: -14706610 -> -9805674
$ws.Cells.Item(74, 14).Value = -9007.75  # N74: -2912.2 -> -9007.75
$ws.Cells.Item(77, 8).Value = 7466419.5  # H77: 11365139 -> 7466419.5
$ws.Cells.Item(77, 9).Value = 9806548  # I77: 14707484 -> 9806548
$ws.Cells.Item(77, 10).Value = 7259.75  # J77: 1164.2 -> 7259.75
$ws.Cells.Item(77, 11).Value = 49032740  # K77: 73537420 -> 49032740
$ws.Cells.Item(77, 12).Value = 36298.75  # L77: 5821 -> 36298.75
$ws.Cells.Item(77, 13).Value = -49028372  # M77: -73533052 -> -49028372
$ws.Cells.Item(77, 14).Value = -45034.75  # N77: -14557 -> -45034.75
$ws.Cells.Item(123, 8).Value = 29678.75  # H123: 29837.428 -> 29678.75
$ws.Cells.Item(123, 10).Value = 29678.75  # J123: 29837.428 -> 29678.75
$ws.Cells.Item(123, 12).Value = 29678.75  # L123: 29837.428 -> 29678.75
$ws.Cells.Item(123, 14).Value = -39478.75  # N123: -39637.428 -> -39478.75
$ws.Cells.Item(125, 8).Value = 32071.2  # H125: 31475.916 -> 32071.2
$ws.Cells.Item(125, 10).Value = 32071.2  # J125: 31475.916 -> 32071.2
$ws.Cells.Item(125, 12).Value = 32071.2  # L125: 31475.916 -> 32071.2
$ws.Cells.Item(125, 14).Value = -41911.2  # N125: -41315.916 -> -41911.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(124, 8).Value = 0  # H124: 30000 -> 0
$ws.Cells.Item(124, 10).Value = 0  # J124: 30000 -> 0
$ws.Cells.Item(124, 12).Value = 0  # L124: 30000 -> 0
$ws.Cells.Item(124, 14).ClearContents()  # N124 removed (was -39820)
$ws.Cells.Item(125, 8).Value = 45781  # H125: 40447 -> 45781
$ws.Cells.Item(125, 10).Value = 45781  # J125: 40447 -> 45781
$ws.Cells.Item(125, 12).Value = 45781  # L125: 40447 -> 45781
$ws.Cells.Item(125, 14).Value = -55621  # N125: -50287 -> -55621

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1881.1234  # H31: 1866.8658 -> 1881.1234
$ws.Cells.Item(31, 9).Value = 1072.7736  # I31: 1066.0927 -> 1072.7736
$ws.Cells.Item(31, 11).Value = 1072.7736  # K31: 1066.0927 -> 1072.7736
$ws.Cells.Item(31, 13).Value = -777.7736  # M31: -771.0926999999999 -> -777.7736
$ws.Cells.Item(34, 8).Value = 1881.1234  # H34: 1866.8658 -> 1881.1234
$ws.Cells.Item(34, 9).Value = 1072.7736  # I34: 1066.0927 -> 1072.7736
$ws.Cells.Item(34, 11).Value = 1072.7736  # K34: 1066.0927 -> 1072.7736
$ws.Cells.Item(34, 13).Value = -870.7736  # M34: -864.0926999999999 -> -870.7736
$ws.Cells.Item(62, 8).Value = 3752.682  # H62: 3942.1667 -> 3752.682
$ws.Cells.Item(62, 9).Value = 2786.5557  # I62: 2854.1428 -> 2786.5557
$ws.Cells.Item(62, 10).Value = 4421.5386  # J62: 4634.5454 -> 4421.5386
$ws.Cells.Item(62, 11).Value = 2786.5557  # K62: 2854.1428 -> 2786.5557
$ws.Cells.Item(62, 12).Value = 4421.5386  # L62: 4634.5454 -> 4421.5386
$ws.Cells.Item(62, 13).Value = -2162.5557  # M62: -2230.1428 -> -2162.5557
$ws.Cells.Item(62, 14).Value = -5669.5386  # N62: -5882.5454 -> -5669.5386
$ws.Cells.Item(65, 8).Value = 3752.682  # H65: 3942.1667 -> 3752.682
$ws.Cells.Item(65, 9).Value = 2786.5557  # I65: 2854.1428 -> 2786.5557
$ws.Cells.Item(65, 10).Value = 4421.5386  # J65: 4634.5454 -> 4421.5386
$ws.Cells.Item(65, 11).Value = 13932.7785  # K65: 14270.714 -> 13932.7785
$ws.Cells.Item(65, 12).Value = 22107.693  # L65: 23172.727 -> 22107.693
$ws.Cells.Item(65, 13).Value = -10812.7785  # M65: -11150.714 -> -10812.7785
$ws.Cells.Item(65, 14).Value = -28347.693  # N65: -29412.727 -> -28347.693
$ws.Cells.Item(122, 8).Value = 50000870  # H122: 53572304 -> 50000870
$ws.Cells.Item(122, 9).Value = 125000610  # I122: 150000580 -> 125000610
$ws.Cells.Item(122, 11).Value = 375001830  # K122: 450001740 -> 375001830
$ws.Cells.Item(122, 13).Value = -374999380  # M122: -449999290 -> -374999380
$ws.Cells.Item(125, 8).Value = 9800  # H125: 0 -> 9800
$ws.Cells.Item(125, 10).Value = 9800  # J125: 0 -> 9800
$ws.Cells.Item(125, 12).Value = 9800  # L125: 0 -> 9800
$ws.Cells.Item(125, 14).Value = -14720  # N125: None -> -14720

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 9649672  # H113: 11458879 -> 9649672
$ws.Cells.Item(113, 9).Value = 6944945  # I113: 8333833 -> 6944945
$ws.Cells.Item(113, 10).Value = 14286348  # J113: 16667289 -> 14286348
$ws.Cells.Item(113, 11).Value = 20834835  # K113: 25001499 -> 20834835
$ws.Cells.Item(113, 12).Value = 42859044  # L113: 50001867 -> 42859044
$ws.Cells.Item(113, 13).Value = -20832665  # M113: -24999329 -> -20832665
$ws.Cells.Item(113, 14).Value = -42863384  # N113: -50006207 -> -42863384
$ws.Cells.Item(131, 8).Value = 692.97  # H131: 684.3200000000001 -> 692.97
$ws.Cells.Item(131, 10).Value = 744.3068  # J131: 734.4773 -> 744.3068
$ws.Cells.Item(131, 12).Value = 2232.9204  # L131: 2203.4319 -> 2232.9204
$ws.Cells.Item(131, 14).Value = -12312.9204  # N131: -12283.4319 -> -12312.9204

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4519.9  # H70: 4133.1665 -> 4519.9
$ws.Cells.Item(70, 9).Value = 4549.8887  # I70: 4118.5625 -> 4549.8887
$ws.Cells.Item(70, 11).Value = 4549.8887  # K70: 4118.5625 -> 4549.8887
$ws.Cells.Item(70, 13).Value = -4279.8887  # M70: -3848.5625 -> -4279.8887
$ws.Cells.Item(73, 8).Value = 4519.9  # H73: 4133.1665 -> 4519.9
$ws.Cells.Item(73, 9).Value = 4549.8887  # I73: 4118.5625 -> 4549.8887
$ws.Cells.Item(73, 11).Value = 4549.8887  # K73: 4118.5625 -> 4549.8887
$ws.Cells.Item(73, 13).Value = -3613.8887  # M73: -3182.5625 -> -3613.8887
$ws.Cells.Item(80, 8).Value = 20005540  # H80: 8336796 -> 20005540
$ws.Cells.Item(80, 9).Value = 8466.333000000001  # I80: 5848.6 -> 8466.333000000001
$ws.Cells.Item(80, 10).Value = 50001150  # J80: 14287472 -> 50001150
$ws.Cells.Item(80, 11).Value = 8466.333000000001  # K80: 5848.6 -> 8466.333000000001
$ws.Cells.Item(80, 12).Value = 50001150  # L80: 14287472 -> 50001150
$ws.Cells.Item(80, 13).Value = -7468.333000000001  # M80: -4850.6 -> -7468.333000000001
$ws.Cells.Item(80, 14).Value = -50003146  # N80: -14289468 -> -50003146
$ws.Cells.Item(83, 8).Value = 20005540  # H83: 8336796 -> 20005540
$ws.Cells.Item(83, 9).Value = 8466.333000000001  # I83: 5848.6 -> 8466.333000000001
$ws.Cells.Item(83, 10).Value = 50001150  # J83: 14287472 -> 50001150
$ws.Cells.Item(83, 11).Value = 42331.665  # K83: 29243 -> 42331.665
$ws.Cells.Item(83, 12).Value = 250005750  # L83: 71437360 -> 250005750
$ws.Cells.Item(83, 13).Value = -37339.665  # M83: -24251 -> -37339.665
$ws.Cells.Item(83, 14).Value = -250015734  # N83: -71447344 -> -250015734
$ws.Cells.Item(102, 8).Value = 919.4400000000001  # H102: 1041.2 -> 919.4400000000001
$ws.Cells.Item(102, 9).Value = 919.4400000000001  # I102: 1041.2 -> 919.4400000000001
$ws.Cells.Item(102, 11).Value = 919.4400000000001  # K102: 1041.2 -> 919.4400000000001
$ws.Cells.Item(102, 13).Value = 702.5599999999999  # M102: 580.8 -> 702.5599999999999
$ws.Cells.Item(124, 8).Value = 0  # H124: 90000 -> 0
$ws.Cells.Item(124, 10).Value = 0  # J124: 90000 -> 0
$ws.Cells.Item(124, 12).Value = 0  # L124: 90000 -> 0
$ws.Cells.Item(124, 14).ClearContents()  # N124 removed (was -99820)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 83334664  # H40: 41668480 -> 83334664
$ws.Cells.Item(40, 9).Value = 1994  # I40: 2178.4 -> 1994
$ws.Cells.Item(40, 11).Value = 1994  # K40: 2178.4 -> 1994
$ws.Cells.Item(40, 13).Value = -1858  # M40: -2042.4 -> -1858
$ws.Cells.Item(122, 8).Value = 8818.9375  # H122: 8409 -> 8818.9375
$ws.Cells.Item(122, 9).Value = 9580.406999999999  # I122: 9298.964 -> 9580.406999999999
$ws.Cells.Item(122, 10).Value = 4707  # J122: 4255.8335 -> 4707
$ws.Cells.Item(122, 11).Value = 28741.221  # K122: 27896.892 -> 28741.221
$ws.Cells.Item(122, 12).Value = 14121  # L122: 12767.5005 -> 14121
$ws.Cells.Item(122, 13).Value = -26291.221  # M122: -25446.892 -> -26291.221
$ws.Cells.Item(122, 14).Value = -19021  # N122: -17667.5005 -> -19021
$ws.Cells.Item(125, 8).Value = 40515.2  # H125: 40715 -> 40515.2
$ws.Cells.Item(125, 10).Value = 40515.2  # J125: 40715 -> 40515.2
$ws.Cells.Item(125, 12).Value = 40515.2  # L125: 40715 -> 40515.2
$ws.Cells.Item(125, 14).Value = -50355.2  # N125: -50555 -> -50355.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 29412216  # H81: 19231218 -> 29412216
$ws.Cells.Item(81, 9).Value = 31250366  # I81: 22727668 -> 31250366
$ws.Cells.Item(81, 10).Value = 1802  # J81: 750.5 -> 1802
$ws.Cells.Item(81, 11).Value = 62500732  # K81: 45455336 -> 62500732
$ws.Cells.Item(81, 12).Value = 3604  # L81: 1501 -> 3604
$ws.Cells.Item(81, 13).Value = -62499671  # M81: -45454275 -> -62499671
$ws.Cells.Item(81, 14).Value = -5726  # N81: -3623 -> -5726
$ws.Cells.Item(84, 8).Value = 29412216  # H84: 19231218 -> 29412216
$ws.Cells.Item(84, 9).Value = 31250366  # I84: 22727668 -> 31250366
$ws.Cells.Item(84, 10).Value = 1802  # J84: 750.5 -> 1802
$ws.Cells.Item(84, 11).Value = 312503660  # K84: 227276680 -> 312503660
$ws.Cells.Item(84, 12).Value = 18020  # L84: 7505 -> 18020
$ws.Cells.Item(84, 13).Value = -312498356  # M84: -227271376 -> -312498356
$ws.Cells.Item(84, 14).Value = -28628  # N84: -18113 -> -28628
$ws.Cells.Item(125, 8).Value = 42153.715  # H125: 43072 -> 42153.715
$ws.Cells.Item(125, 10).Value = 42153.715  # J125: 43072 -> 42153.715
$ws.Cells.Item(125, 12).Value = 42153.715  # L125: 43072 -> 42153.715
$ws.Cells.Item(125, 14).Value = -51993.715  # N125: -52912 -> -51993.715
